# Update "views/likes" counts in column F on the "展览" (exhibitions) sheet
# and the "全部类型" (all types) sheet, per the upstream gh-pages data refresh.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$wsExhibit.Range("F2").Value  = 131
$wsExhibit.Range("F3").Value  = 233
$wsExhibit.Range("F4").Value  = 14
$wsExhibit.Range("F5").Value  = 6683
$wsExhibit.Range("F8").Value  = 137
$wsExhibit.Range("F9").Value  = 6172
$wsExhibit.Range("F12").Value = 1251
$wsExhibit.Range("F13").Value = 1251
$wsExhibit.Range("F15").Value = 95
$wsExhibit.Range("F17").Value = 119
$wsExhibit.Range("F19").Value = 361
$wsExhibit.Range("F21").Value = 7
$wsExhibit.Range("F22").Value = 4500
$wsExhibit.Range("F23").Value = 53
$wsExhibit.Range("F24").Value = 30
$wsExhibit.Range("F25").Value = 190
$wsExhibit.Range("F26").Value = 50

# --- 全部类型 (sheet4) ---
$wsAll.Range("F2").Value  = 131
$wsAll.Range("F3").Value  = 233
$wsAll.Range("F4").Value  = 14
$wsAll.Range("F5").Value  = 6683
$wsAll.Range("F8").Value  = 137
$wsAll.Range("F9").Value  = 6172
$wsAll.Range("F12").Value = 1251
$wsAll.Range("F13").Value = 1251
$wsAll.Range("F15").Value = 95
$wsAll.Range("F17").Value = 119
$wsAll.Range("F19").Value = 361
$wsAll.Range("F21").Value = 7
$wsAll.Range("F22").Value = 4501
$wsAll.Range("F24").Value = 53
$wsAll.Range("F25").Value = 30
$wsAll.Range("F26").Value = 190
$wsAll.Range("F27").Value = 50
